# quarterly.xlsx update: roll the quarterly window forward by two quarters
# (drop the two oldest quarters -> 1399/04, 1399/07; append the two newest
# quarters -> 1401/10, 1402/01) and refresh the figures for every quarter
# that is still in range ("update database and change read_price algorithm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header rows (8 and 24): the 10 quarter-headers in E:N shift two columns
#    to the left (E<-G, F<-H, ... L<-N) and two brand-new quarter labels are
#    appended in M and N. Derive the two brand-new labels from the existing
#    labels via text substitution so the Persian text is taken verbatim from
#    the workbook itself rather than retyped.
# ---------------------------------------------------------------------------
$headerRows = @(8, 24)

foreach ($r in $headerRows) {
    # Capture the current (pre-shift) header text for columns E..N first.
    $oldE = $ws.Cells.Item($r, 5).Value2   # E - "... 1399/04"
    $oldF = $ws.Cells.Item($r, 6).Value2   # F - "... 1399/07"
    $oldG = $ws.Cells.Item($r, 7).Value2   # G - "... 1399/10"
    $oldH = $ws.Cells.Item($r, 8).Value2   # H - "... 1400/01"
    $oldI = $ws.Cells.Item($r, 9).Value2   # I - "... 1400/04"
    $oldJ = $ws.Cells.Item($r, 10).Value2  # J - "... 1400/07"
    $oldK = $ws.Cells.Item($r, 11).Value2  # K - "... 1400/10"
    $oldL = $ws.Cells.Item($r, 12).Value2  # L - "... 1401/01"
    $oldM = $ws.Cells.Item($r, 13).Value2  # M - "... 1401/04"
    $oldN = $ws.Cells.Item($r, 14).Value2  # N - "... 1401/07"

    # New quarter labels for the two newly-added columns, built from the
    # existing labels so the Persian wording/formatting matches exactly.
    $newM = $oldG.Replace("1399/10", "1401/10")  # "فصل چهارم منتهی به 1401/10"
    $newN = $oldH.Replace("1400/01", "1402/01")  # "فصل اول منتهی به 1402/01"

    $ws.Cells.Item($r, 5).Value  = $oldG   # E <- old G
    $ws.Cells.Item($r, 6).Value  = $oldH   # F <- old H
    $ws.Cells.Item($r, 7).Value  = $oldI   # G <- old I
    $ws.Cells.Item($r, 8).Value  = $oldJ   # H <- old J
    $ws.Cells.Item($r, 9).Value  = $oldK   # I <- old K
    $ws.Cells.Item($r, 10).Value = $oldL   # J <- old L
    $ws.Cells.Item($r, 11).Value = $oldM   # K <- old M
    $ws.Cells.Item($r, 12).Value = $oldN   # L <- old N
    $ws.Cells.Item($r, 13).Value = $newM   # M <- new quarter (1401/10)
    $ws.Cells.Item($r, 14).Value = $newN   # N <- new quarter (1402/01)
}

# ---------------------------------------------------------------------------
# 2) Data rows: refreshed quarterly figures (database update + the new
#    read_price algorithm recompute some historical quarters, not just a
#    plain left-shift), row numbers map to the same rows as the headers.
# ---------------------------------------------------------------------------
$data = @{}
$data[10] = @(159739,146969,61741,103682,228365,231434,433894,330285,89366,151626)
$data[11] = @(0,0,0,0,0,0,0,0,0,0)
$data[12] = @(0,0,0,0,0,0,0,0,0,0)
$data[13] = @(-337,0,0,0,0,0,0,0,0,0)
$data[14] = @(0,0,0,0,0,0,0,0,0,0)
$data[15] = @(58,0,151,86,-237,0,0,0,0,92)
$data[16] = @(242,613,8991,3227,-9313,125,1006,515,-1240,704)
$data[17] = @(19299,24921,4394,33638,11123,20149,34313,35618,31279,31839)
$data[18] = @(6092,0,0,0,25788,0,0,0,1119,0)
$data[19] = @(54866,30728,83573,22768,-13015,51145,34574,41567,35849,41539)
$data[20] = @(239959,203231,158850,163401,242711,302853,503787,407985,156373,225800)
$data[26] = @(58,61,61,61,67,67,51,51,65,61)
$data[27] = @(386,349,374,374,674,374,377,377,375,367)

$cols = @(5,6,7,8,9,10,11,12,13,14)  # E..N

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $vals[$i]
    }
}
